$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: row 3 is the b.md entry. Mark it as handed off. ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-28 10:45:07"

# --- "zh-cn" sheet: row 3 is the b.md entry. ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 10:44:57"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4dd47b1923181cbd050d11c94e46ba5693cfa402/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/429075ae3e35b60f70a560b5c683278c1253c445/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- "de-de" sheet: row 3 is the b.md entry. ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-28 10:45:07"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4dd47b1923181cbd050d11c94e46ba5693cfa402/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/429075ae3e35b60f70a560b5c683278c1253c445/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
